$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1152.0435
$ws.Range("I17").Value = 99.75
$ws.Range("J17").Value = 1373.579
$ws.Range("K17").Value = 299.25
$ws.Range("L17").Value = 4120.737
$ws.Range("M17").Value = -131.25
$ws.Range("N17").Value = -4456.737
# Row 21
$ws.Range("H21").Value = 3633
$ws.Range("I21").Value = 3633
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 3633
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -3165
# Row 23
$ws.Range("H23").Value = 3633
$ws.Range("I23").Value = 3633
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 3633
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3399
# Row 34
$ws.Range("H34").Value = 7490.4
$ws.Range("I34").Value = 6724.125
$ws.Range("J34").Value = 10555.5
$ws.Range("K34").Value = 6724.125
$ws.Range("L34").Value = 10555.5
$ws.Range("M34").Value = -6521.125
$ws.Range("N34").Value = -10961.5
# Row 36
$ws.Range("H36").Value = 7490.4
$ws.Range("I36").Value = 6724.125
$ws.Range("J36").Value = 10555.5
$ws.Range("K36").Value = 6724.125
$ws.Range("L36").Value = 10555.5
$ws.Range("M36").Value = -6009.125
$ws.Range("N36").Value = -11985.5
# Row 61
$ws.Range("H61").Value = 4005
$ws.Range("I61").Value = 4005
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 12015
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -11843
# Row 137
$ws.Range("H137").Value = 2356
$ws.Range("I137").Value = 1716.3846
$ws.Range("J137").Value = 3048.9167
$ws.Range("K137").Value = 5149.1538
$ws.Range("L137").Value = 9146.750100000001
$ws.Range("M137").Value = -2599.1538
$ws.Range("N137").Value = -14246.7501
# Row 138
$ws.Range("H138").Value = 3230052
$ws.Range("I138").Value = 2876.8462
$ws.Range("J138").Value = 5560789.5
$ws.Range("K138").Value = 8630.5386
$ws.Range("L138").Value = 16682368.5
$ws.Range("M138").Value = -3490.5386
$ws.Range("N138").Value = -16692648.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1869155.9
$ws.Range("I32").Value = 2407108
$ws.Range("J32").Value = 4254.6665
$ws.Range("K32").Value = 2407108
$ws.Range("L32").Value = 4254.6665
$ws.Range("M32").Value = -2406821
$ws.Range("N32").Value = -4828.6665
# Row 61
$ws.Range("H61").Value = 5902.018
$ws.Range("I61").Value = 2750.853
$ws.Range("J61").Value = 11003.904
$ws.Range("K61").Value = 2750.853
$ws.Range("L61").Value = 11003.904
$ws.Range("M61").Value = -2538.853
$ws.Range("N61").Value = -11427.904
# Row 110
$ws.Range("H110").Value = 37038384
$ws.Range("I110").Value = 1239
$ws.Range("J110").Value = 66668100
$ws.Range("K110").Value = 1239
$ws.Range("L110").Value = 66668100
$ws.Range("M110").Value = 806
$ws.Range("N110").Value = -66672190
# Row 132
$ws.Range("H132").Value = 1932044.1
$ws.Range("I132").Value = 7156473
$ws.Range("J132").Value = 7254.5264
$ws.Range("K132").Value = 21469419
$ws.Range("L132").Value = 21763.5792
$ws.Range("M132").Value = -21466889
$ws.Range("N132").Value = -26823.5792
# Row 136
$ws.Range("H136").Value = 5902.018
$ws.Range("I136").Value = 2750.853
$ws.Range("J136").Value = 11003.904
$ws.Range("K136").Value = 8252.559000000001
$ws.Range("L136").Value = 33011.712
$ws.Range("M136").Value = -5702.559000000001
$ws.Range("N136").Value = -38111.712

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 7580607
$ws.Range("I99").Value = 4366.8335
$ws.Range("J99").Value = 15156848
$ws.Range("K99").Value = 4366.8335
$ws.Range("L99").Value = 15156848
$ws.Range("M99").Value = -2868.8335
$ws.Range("N99").Value = -15159844
# Row 132
$ws.Range("H132").Value = 178852.33
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 178852.33
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 178852.33
$ws.Range("N132").Value = -188972.33
# Row 134
$ws.Range("H134").Value = 7915.8096
$ws.Range("I134").Value = 3362.4167
$ws.Range("J134").Value = 13987
$ws.Range("K134").Value = 10087.2501
$ws.Range("L134").Value = 41961
$ws.Range("M134").Value = -7552.250100000001
$ws.Range("N134").Value = -47031

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 2514.2856
$ws.Range("I58").Value = 2514.2856
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 7542.8568
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -7414.8568
# Row 68
$ws.Range("H68").Value = 2480.2307
$ws.Range("I68").Value = 2399.8
$ws.Range("J68").Value = 2499.3809
$ws.Range("K68").Value = 7199.400000000001
$ws.Range("L68").Value = 7498.1427
$ws.Range("M68").Value = -6388.400000000001
$ws.Range("N68").Value = -9120.1427
# Row 71
$ws.Range("H71").Value = 2480.2307
$ws.Range("I71").Value = 2399.8
$ws.Range("J71").Value = 2499.3809
$ws.Range("K71").Value = 21598.2
$ws.Range("L71").Value = 22494.4281
$ws.Range("M71").Value = -17542.2
$ws.Range("N71").Value = -30606.4281
# Row 107
$ws.Range("H107").Value = 4652635.5
$ws.Range("I107").Value = 783.3333
$ws.Range("J107").Value = 5884008.5
$ws.Range("K107").Value = 2349.9999
$ws.Range("L107").Value = 17652025.5
$ws.Range("M107").Value = -429.9998999999998
$ws.Range("N107").Value = -17655865.5
# Row 122
$ws.Range("H122").Value = 944138.0600000001
$ws.Range("I122").Value = 2572651
$ws.Range("J122").Value = 1314.7368
$ws.Range("K122").Value = 23153859
$ws.Range("L122").Value = 11832.6312
$ws.Range("M122").Value = -23151409
$ws.Range("N122").Value = -16732.6312

$ws = $wb.Worksheets.Item("GSM")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 102
$ws.Range("H102").Value = 8652.143
$ws.Range("I102").Value = 8427.5
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 8427.5
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -6805.5
$ws.Range("N102").Value = -13244
# Row 104
$ws.Range("H104").Value = 47125
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 47125
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 47125
$ws.Range("N104").Value = -54113
# Row 126
$ws.Range("H126").Value = 5870.273
$ws.Range("I126").Value = 2893.25
$ws.Range("J126").Value = 7571.4287
$ws.Range("K126").Value = 8679.75
$ws.Range("L126").Value = 22714.2861
$ws.Range("M126").Value = -6209.75
$ws.Range("N126").Value = -27654.2861
# Row 140
$ws.Range("H140").Value = 44997
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 44997
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 44997
$ws.Range("N140").Value = -55357

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2526.2
$ws.Range("I16").Value = 2610.9443
$ws.Range("J16").Value = 1763.5
$ws.Range("K16").Value = 2610.9443
$ws.Range("L16").Value = 1763.5
$ws.Range("M16").Value = -2440.9443
$ws.Range("N16").Value = -2103.5
# Row 100
$ws.Range("H100").Value = 5333.3335
$ws.Range("I100").Value = 3249.25
$ws.Range("J100").Value = 7000.6
$ws.Range("K100").Value = 3249.25
$ws.Range("L100").Value = 7000.6
$ws.Range("M100").Value = -2708.25
$ws.Range("N100").Value = -8082.6
# Row 123
$ws.Range("H123").Value = 51528
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 51528
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 51528
$ws.Range("N123").Value = -61328
# Row 139
$ws.Range("H139").Value = 64515
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 64515
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 64515
$ws.Range("N139").Value = -74795
# Row 141
$ws.Range("H141").Value = 75179.164
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 75179.164
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 75179.164
$ws.Range("N141").Value = -85539.164

$ws = $wb.Worksheets.Item("WVR")
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
# Row 133
$ws.Range("H133").Value = 134983
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 134983
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 134983
$ws.Range("N133").Value = -145103
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 139
$ws.Range("H139").Value = 68143.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 68143.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 68143.336
$ws.Range("N139").Value = -78423.336
# Row 140
$ws.Range("H140").Value = 99981.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99981.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99981.5
$ws.Range("N140").Value = -110341.5
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

